$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MM2300400043"
$ws.Range("A3").Value = "MM2300400047"
$ws.Range("A4").Value = "MM2300400048"
$ws.Range("A5").Value = "MM2300400049"
$ws.Range("A6").Value = "MM2316800128"

$ws.Cells.Item(11, 8).Select()
